$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "71.100.62"
$ws.Range("E2").Value = "  +5.97%  "

# Row 3
$ws.Range("D3").Value = "3.713.06"
$ws.Range("E3").Value = "  +19.70%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'620.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.21%  "

# Row 6
$ws.Range("D6").Value = "'183.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.19%  "

# Row 7
$ws.Range("D7").Value = "3.707.07"
$ws.Range("E7").Value = "  +19.59%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").Value = "'0.543"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.70%  "

# Row 10
$ws.Range("E10").Value = "  +8.39%  "

# Row 11
$ws.Range("E11").Value = "  +4.18%  "

# Row 12
$ws.Range("D12").Value = "'0.504"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.69%  "

# Row 13
$ws.Range("D13").Value = "'40.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +12.93%  "

# Row 14
$ws.Range("E14").Value = "  +6.37%  "

# Row 15
$ws.Range("D15").Value = "4.332.03"
$ws.Range("E15").Value = "  +19.46%  "

# Row 16
$ws.Range("D16").Value = "3.711.92"
$ws.Range("E16").Value = "  +19.48%  "

# Row 17
$ws.Range("D17").Value = "71.159.43"
$ws.Range("E17").Value = "  +6.08%  "

# Row 18
$ws.Range("E18").Value = "  +1.79%  "

# Row 19
$ws.Range("D19").Value = "'7.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.22%  "

# Row 20
$ws.Range("D20").Value = "'518.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.72%  "

# Row 21
$ws.Range("D21").Value = "'16.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.68%  "

# Row 22
$ws.Range("D22").Value = "'9.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +19.98%  "

# Row 23
$ws.Range("E23").Value = "  +8.46%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'88.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.29%  "

# Row 25
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D25").Value = "'2.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.04%  "

# Row 26
$ws.Range("D26").Value = "'13.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.13%  "

# Row 27
$ws.Range("D27").Value = "'11.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.28%  "

# Row 28
$ws.Range("E28").Value = "  +0.09%  "

# Row 29
$ws.Range("E29").Value = "  +10.99%  "

# Row 30
$ws.Range("D30").Value = "'8.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.77%  "

# Row 31
$ws.Range("E31").Value = "  +11.99%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'31.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.15%  "

# Row 33
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "'0.0000111"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +18.38%  "

# Row 34
$ws.Range("E34").Value = "  +4.57%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("E36").Value = "  +10.10%  "

# Row 37
$ws.Range("E37").Value = "  +9.87%  "

# Row 38
$ws.Range("D38").Value = "'0.347"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.72%  "

# Row 39
$ws.Range("E39").Value = "  +12.49%  "

# Row 40
$ws.Range("D40").Value = "'0.135"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.02%  "

# Row 41
$ws.Range("D41").Value = "'51.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.39%  "

# Row 42
$ws.Range("D42").Value = "'436.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +17.76%  "

# Row 43
$ws.Range("D43").Value = "'45.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.48%  "

# Row 44
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").Value = "'8.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.87%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "3.144.73"
$ws.Range("E45").Value = "  +12.65%  "

# Row 46
$ws.Range("D46").Value = "'2.87"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.24%  "

# Row 47
$ws.Range("D47").Value = "'0.0370"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.42%  "

# Row 48
$ws.Range("D48").Value = "'28.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +11.26%  "

# Row 49
$ws.Range("D49").Value = "'140.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.79%  "

# Row 51
$ws.Range("E51").Value = "  +8.44%  "
